$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1878612716763006
$ws.Range("C2").Value = 0.5578034682080925
$ws.Range("J2").Value = 0.005780346820809248
$ws.Range("P2").Value = 0.1445086705202312
$ws.Range("S2").Value = 0.1040462427745665
$ws.Range("B3").Value = 0.0198019801980198
$ws.Range("C3").Value = 0.04455445544554455
$ws.Range("J3").Value = 0.01485148514851485
$ws.Range("P3").Value = 0.7079207920792079
$ws.Range("S3").Value = 0.2128712871287129
$ws.Range("P4").Value = 0.8378378378378378
$ws.Range("S4").Value = 0.1621621621621622
$ws.Range("J5").Value = 0.2
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.04504504504504504
$ws.Range("D6").Value = 0.009009009009009009
$ws.Range("E6").Value = 0.009009009009009009
$ws.Range("F6").Value = 0.04954954954954955
$ws.Range("J6").Value = 0.2477477477477477
$ws.Range("O6").Value = 0.03603603603603604
$ws.Range("Q6").Value = 0.1576576576576577
$ws.Range("R6").Value = 0.08558558558558559
$ws.Range("S6").Value = 0.3603603603603603
$ws.Range("B7").Value = 0.1261682242990654
$ws.Range("D7").Value = 0.009345794392523364
$ws.Range("E7").Value = 0.004672897196261682
$ws.Range("F7").Value = 0.07009345794392523
$ws.Range("J7").Value = 0.1588785046728972
$ws.Range("O7").Value = 0.01869158878504673
$ws.Range("Q7").Value = 0.1635514018691589
$ws.Range("R7").Value = 0.1074766355140187
$ws.Range("S7").Value = 0.3411214953271028
$ws.Range("B8").Value = 0.1061151079136691
$ws.Range("D8").Value = 0.01618705035971223
$ws.Range("E8").Value = 0.001798561151079137
$ws.Range("F8").Value = 0.05935251798561151
$ws.Range("J8").Value = 0.1366906474820144
$ws.Range("O8").Value = 0.0197841726618705
$ws.Range("Q8").Value = 0.1528776978417266
$ws.Range("R8").Value = 0.08453237410071943
$ws.Range("S8").Value = 0.4226618705035971
$ws.Range("B9").Value = 0.1228070175438596
$ws.Range("D9").Value = 0.03508771929824561
$ws.Range("F9").Value = 0.06140350877192982
$ws.Range("J9").Value = 0.1271929824561404
$ws.Range("O9").Value = 0.0131578947368421
$ws.Range("Q9").Value = 0.1578947368421053
$ws.Range("R9").Value = 0.07017543859649122
$ws.Range("S9").Value = 0.412280701754386
$ws.Range("B10").Value = 0.1151750972762646
$ws.Range("D10").Value = 0.01478599221789883
$ws.Range("E10").Value = 0.002334630350194552
$ws.Range("F10").Value = 0.07392996108949416
$ws.Range("J10").Value = 0.1346303501945525
$ws.Range("O10").Value = 0.01712062256809339
$ws.Range("Q10").Value = 0.1828793774319066
$ws.Range("R10").Value = 0.09571984435797666
$ws.Range("S10").Value = 0.3634241245136187
$ws.Range("G11").Value = 0.1411042944785276
$ws.Range("J11").Value = 0.1012269938650307
$ws.Range("K11").Value = 0.2269938650306748
$ws.Range("L11").Value = 0.5153374233128835
$ws.Range("S11").Value = 0.01533742331288344
$ws.Range("G12").Value = 0.7527472527472527
$ws.Range("J12").Value = 0.1593406593406593
$ws.Range("K12").Value = 0.005494505494505495
$ws.Range("L12").Value = 0.03296703296703297
$ws.Range("S12").Value = 0.04945054945054945
$ws.Range("G13").Value = 0.7346938775510204
$ws.Range("J13").Value = 0.1836734693877551
$ws.Range("S13").Value = 0.08163265306122448
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.03043478260869565
$ws.Range("H15").Value = 0.1478260869565217
$ws.Range("I15").Value = 0.08260869565217391
$ws.Range("J15").Value = 0.3478260869565217
$ws.Range("K15").Value = 0.05217391304347826
$ws.Range("M15").Value = 0.01739130434782609
$ws.Range("O15").Value = 0.04347826086956522
$ws.Range("S15").Value = 0.2782608695652174
$ws.Range("F16").Value = 0.004672897196261682
$ws.Range("H16").Value = 0.2570093457943925
$ws.Range("I16").Value = 0.1074766355140187
$ws.Range("J16").Value = 0.3224299065420561
$ws.Range("K16").Value = 0.1121495327102804
$ws.Range("M16").Value = 0.02803738317757009
$ws.Range("O16").Value = 0.04205607476635514
$ws.Range("S16").Value = 0.1261682242990654
$ws.Range("F17").Value = 0.009324009324009324
$ws.Range("H17").Value = 0.2121212121212121
$ws.Range("I17").Value = 0.09324009324009325
$ws.Range("J17").Value = 0.37995337995338
$ws.Range("K17").Value = 0.1212121212121212
$ws.Range("M17").Value = 0.02331002331002331
$ws.Range("N17").Value = 0.002331002331002331
$ws.Range("O17").Value = 0.05361305361305362
$ws.Range("S17").Value = 0.1048951048951049
$ws.Range("F18").Value = 0.01731601731601732
$ws.Range("H18").Value = 0.2164502164502164
$ws.Range("I18").Value = 0.09090909090909091
$ws.Range("J18").Value = 0.3463203463203463
$ws.Range("K18").Value = 0.1168831168831169
$ws.Range("M18").Value = 0.008658008658008658
$ws.Range("O18").Value = 0.06060606060606061
$ws.Range("S18").Value = 0.1428571428571428
$ws.Range("F19").Value = 0.01080691642651297
$ws.Range("H19").Value = 0.2384726224783862
$ws.Range("I19").Value = 0.0893371757925072
$ws.Range("J19").Value = 0.3371757925072046
$ws.Range("K19").Value = 0.09293948126801153
$ws.Range("M19").Value = 0.01945244956772334
$ws.Range("O19").Value = 0.07060518731988473
$ws.Range("S19").Value = 0.1412103746397695
